# MIC1 microinstruction synthesizer - fill in a new microinstruction
# (NEXT_ADDR hex = 38, ALU op = A AND B, MEM = read, C bus select = CPP,
#  and mark the extra JAM bit in Q3), and drop the sheet's password
# protection (keeping the sheet otherwise as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships password-protected; several of the cells we need to
# touch (the G2 ALU-select formula in particular) are locked, so drop the
# protection before editing. This also matches the author's removal of the
# stored password from the sheet.
$ws.Unprotect()

# NEXT_ADDR (hex, stored as a plain number - HEX2BIN treats the digits as hex)
$ws.Range("A2").Value = 38

# Extend the ALU dropdown formula with the new "A AND B" -> 12 branch.
$ws.Range("G2").Formula = '=IF(H2="A",24,IF(H2="B",20,IF(H2="AN",26,IF(H2="BN",44,IF(H2="A+B",60,IF(H2="B+1",53,IF(H2="0",16,IF(H2="B-1",54,IF(H2="A OR B",28,IF(H2="B-A",63,IF(H2="A AND B",12,0)))))))))))'

# Fill in the MEM / ALU / C-bus-select inputs. "read" must be written before
# "A AND B" so the shared-string table order matches the saved workbook.
$ws.Range("R2").Value = "read"
$ws.Range("H2").Value = "A AND B"
$ws.Range("U2").Value = "CPP"

# JAM bit (Q3) flips on for this microinstruction.
$ws.Range("Q3").Value = 1

# Move the active selection the way the author left it.
$ws.Range("G3:H3").Select()

$wb.Save()
